$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue $ws 'D2' '63.529.40'
$ws.Range('E2').Value = '  +2.80%  '
Set-TextValue $ws 'D3' '3.475.73'
$ws.Range('E3').Value = '  +1.72%  '
Set-TextValue $ws 'D4' '0.999'
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue $ws 'D5' '581.62'
$ws.Range('E5').Value = '  +0.59%  '
Set-TextValue $ws 'D6' '147.53'
$ws.Range('E6').Value = '  +1.68%  '
Set-TextValue $ws 'D7' '3.474.98'
$ws.Range('E7').Value = '  +1.65%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('E12').Value = '  +4.93%  '
Set-TextValue $ws 'D13' '4.072.01'
$ws.Range('E13').Value = '  +1.79%  '
Set-TextValue $ws 'D14' '29.71'
$ws.Range('E14').Value = '  +5.51%  '
$ws.Range('E15').Value = '  +2.36%  '
Set-TextValue $ws 'D16' '3.481.41'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('E17').Value = '  +1.39%  '
Set-TextValue $ws 'D18' '63.448.01'
$ws.Range('E18').Value = '  +2.65%  '
Set-TextValue $ws 'D19' '6.35'
$ws.Range('E19').Value = '  +3.07%  '
$ws.Range('E20').Value = '  +3.78%  '
Set-TextValue $ws 'D21' '9.35'
$ws.Range('E21').Value = '  +1.87%  '
Set-TextValue $ws 'D22' '390.09'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E23').Value = '  +2.52%  '
Set-TextValue $ws 'D24' '75.21'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('E25').Value = '  -0.08%  '
Set-TextValue $ws 'D26' '3.623.44'
$ws.Range('E26').Value = '  +1.92%  '
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('E28').Value = '  -4.36%  '
$ws.Range('E29').Value = '  +2.24%  '
Set-TextValue $ws 'D30' '1.00'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  +3.05%  '
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  -3.81%  '
Set-TextValue $ws 'D35' '23.56'
$ws.Range('E35').Value = '  +0.31%  '
Set-TextValue $ws 'B36' 'NEARProtocol'
Set-TextValue $ws 'C36' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D36' '5.32'
$ws.Range('E36').Value = '  +1.57%  '
Set-TextValue $ws 'B37' 'Aptos'
Set-TextValue $ws 'C37' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D37' '7.13'
$ws.Range('E37').Value = '  +2.33%  '
$ws.Range('E38').Value = '  +8.49%  '
Set-TextValue $ws 'D39' '31.67'
$ws.Range('E39').Value = '  +9.89%  '
Set-TextValue $ws 'D40' '169.72'
$ws.Range('E40').Value = '  +0.64%  '
Set-TextValue $ws 'D41' '3.513.62'
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('E43').Value = '  +1.63%  '
Set-TextValue $ws 'D44' '1.73'
$ws.Range('E44').Value = '  +3.69%  '
Set-TextValue $ws 'D45' '42.41'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('E46').Value = '  +3.38%  '
$ws.Range('E47').Value = '  -0.58%  '
Set-TextValue $ws 'D48' '2.613.11'
$ws.Range('E48').Value = '  +4.16%  '
$ws.Range('E49').Value = '  +9.50%  '
Set-TextValue $ws 'D50' '23.13'
$ws.Range('E50').Value = '  +1.51%  '
Set-TextValue $ws 'D51' '6.78'
$ws.Range('E51').Value = '  +2.47%  '
